$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row and municipality/state name casing ("de"/"del"/"el"/"la"/"los"/"las" -> capitalized)
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'
$ws.Range('B4').Value = 'Rincón De Romos'
$ws.Range('B19').Value = 'Amatenango De La Frontera'
$ws.Range('B28').Value = 'Comitán De Domínguez'
$ws.Range('B49').Value = 'Ocozocoautla De Espinosa'
$ws.Range('B54').Value = 'Salto De Agua'
$ws.Range('B55').Value = 'San Cristóbal De Las Casas'
$ws.Range('B81').Value = 'Hidalgo Del Parral'
$ws.Range('B85').Value = 'San Francisco De Borja'
$ws.Range('B99').Value = 'San Juan De Sabinas'
$ws.Range('A107').Value = 'Ciudad De México'
$ws.Range('B126').Value = 'Pánuco De Coronado'
$ws.Range('B128').Value = 'San Juan De Guadalupe'
$ws.Range('A133').Value = 'Estado De México'
$ws.Range('B133').Value = 'Almoloya De Alquisiras'
$ws.Range('B141').Value = 'Ecatepec De Morelos'
$ws.Range('B152').Value = 'Naucalpan De Juárez'
$ws.Range('B156').Value = 'San Felipe Del Progreso'
$ws.Range('B165').Value = 'Tlalnepantla De Baz'
$ws.Range('B170').Value = 'Valle De Bravo'
$ws.Range('B171').Value = 'Villa Del Carbón'
$ws.Range('B178').Value = 'San Miguel De Allende'
$ws.Range('B179').Value = 'Apaseo El Alto'
$ws.Range('B180').Value = 'Apaseo El Grande'
$ws.Range('B187').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B197').Value = 'San Diego De La Unión'
$ws.Range('B199').Value = 'San Francisco Del Rincón'
$ws.Range('B201').Value = 'San Luis De La Paz'
$ws.Range('B202').Value = 'Santa Cruz De Juventino Rosas'
$ws.Range('B203').Value = 'Silao De La Victoria'
$ws.Range('B207').Value = 'Valle De Santiago'
$ws.Range('B213').Value = 'Acapulco De Juárez'
$ws.Range('B216').Value = 'Ajuchitlán Del Progreso'
$ws.Range('B217').Value = 'Alcozauca De Guerrero'
$ws.Range('B220').Value = 'Atoyac De Álvarez'
$ws.Range('B221').Value = 'Ayutla De Los Libres'
$ws.Range('B223').Value = 'Chilapa De Álvarez'
$ws.Range('B224').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B227').Value = 'Coyuca De Benítez'
$ws.Range('B228').Value = 'Cuetzala Del Progreso'
$ws.Range('B230').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B231').Value = 'Iguala De La Independencia'
$ws.Range('B235').Value = 'Mártir De Cuilapan'
$ws.Range('B239').Value = 'Taxco De Alarcón'
$ws.Range('B240').Value = 'Técpan De Galeana'
$ws.Range('B241').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B243').Value = 'Tixtla De Guerrero'
$ws.Range('B247').Value = 'Tlapa De Comonfort'
$ws.Range('B255').Value = 'Agua Blanca De Iturbide'
$ws.Range('B267').Value = 'Jacala De Ledezma'
$ws.Range('B272').Value = 'Mineral Del Chico'
$ws.Range('B273').Value = 'Mineral Del Monte'
$ws.Range('B274').Value = 'Mixquiahuala De Juárez'
$ws.Range('B275').Value = 'Molango De Escamilla'
$ws.Range('B277').Value = 'Pachuca De Soto'
$ws.Range('B281').Value = 'Santiago De Anaya'
$ws.Range('B283').Value = 'Tenango De Doria'
$ws.Range('B285').Value = 'Tepeji Del Río De Ocampo'
$ws.Range('B289').Value = 'Tula De Allende'
$ws.Range('B290').Value = 'Tulancingo De Bravo'
$ws.Range('B292').Value = 'Zacualtipán De Ángeles'
$ws.Range('B296').Value = 'Atotonilco El Alto'
$ws.Range('B301').Value = 'Encarnación De Díaz'
$ws.Range('B305').Value = 'Huejuquilla El Alto'
$ws.Range('B306').Value = 'Ixtlahuacán De Los Membrillos'
$ws.Range('B308').Value = 'Jilotlán De Los Dolores'
$ws.Range('B310').Value = 'Lagos De Moreno'
$ws.Range('B315').Value = 'Tamazula De Gordiano'
$ws.Range('B318').Value = 'Teocuitatlán De Corona'
$ws.Range('B319').Value = 'Tepatitlán De Morelos'
$ws.Range('B321').Value = 'Tizapán El Alto'
$ws.Range('B326').Value = 'Unión De San Antonio'
$ws.Range('B329').Value = 'Zapotlán El Grande'
$ws.Range('B337').Value = 'Coalcomán De Vázquez Pallares'
$ws.Range('B384').Value = 'Puente De Ixtla'
$ws.Range('B393').Value = 'Ixtlán Del Río'
$ws.Range('B409').Value = 'San Nicolás De Los Garza'
$ws.Range('B411').Value = 'Acatlán De Pérez Figueroa'
$ws.Range('B414').Value = 'Coicoyán De Las Flores'
$ws.Range('B416').Value = 'Constancia Del Rosario'
$ws.Range('B418').Value = 'Ixtlán De Juárez'
$ws.Range('B419').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B422').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B423').Value = 'Oaxaca De Juárez'
$ws.Range('B424').Value = 'Putla Villa De Guerrero'
$ws.Range('B453').Value = 'Santa Lucía Del Camino'
$ws.Range('B470').Value = 'Tlacolula De Matamoros'
$ws.Range('B471').Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range('B489').Value = 'Izúcar De Matamoros'
$ws.Range('B493').Value = 'Palmar De Bravo'
$ws.Range('B500').Value = 'Tecali De Herrera'
$ws.Range('B502').Value = 'Tepexi De Rodríguez'
$ws.Range('B503').Value = 'Tetela De Ocampo'
$ws.Range('B518').Value = 'Amealco De Bonfil'
$ws.Range('B520').Value = 'Cadereyta De Montes'
$ws.Range('B523').Value = 'Jalpan De Serra'
$ws.Range('B524').Value = 'Landa De Matamoros'
$ws.Range('B526').Value = 'Pinal De Amoles'
$ws.Range('B529').Value = 'San Juan Del Río'
$ws.Range('B539').Value = 'Cerro De San Pedro'
$ws.Range('B540').Value = 'Ciudad Del Maíz'
$ws.Range('B550').Value = 'San Ciro De Acosta'
$ws.Range('B554').Value = 'Santa María Del Río'
$ws.Range('B561').Value = 'Tanquián De Escobedo'
$ws.Range('B564').Value = 'Villa De Arriaga'
$ws.Range('B565').Value = 'Villa De Ramos'
$ws.Range('B566').Value = 'Villa De Reyes'
$ws.Range('B580').Value = 'Nacozari De García'
$ws.Range('B590').Value = 'Jalpa De Méndez'
$ws.Range('B616').Value = 'Soto La Marina'
$ws.Range('B627').Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range('B634').Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range('B638').Value = 'Amatlán De Los Reyes'
$ws.Range('B645').Value = 'Castillo De Teayo'
$ws.Range('B652').Value = 'Cosamaloapan De Carpio'
$ws.Range('B653').Value = 'Cosautlán De Carvajal'
$ws.Range('B663').Value = 'Hueyapan De Ocampo'
$ws.Range('B665').Value = 'Ixhuatlán De Madero'
$ws.Range('B674').Value = 'Lerdo De Tejada'
$ws.Range('B675').Value = 'Martínez De La Torre'
$ws.Range('B677').Value = 'Medellín De Bravo'
$ws.Range('B680').Value = 'Mixtla De Altamirano'
$ws.Range('B687').Value = 'Paso Del Macho'
$ws.Range('B691').Value = 'Poza Rica De Hidalgo'
$ws.Range('B697').Value = 'Sayula De Alemán'
$ws.Range('B698').Value = 'Soledad De Doblado'
$ws.Range('B717').Value = 'Vega De Alatorre'
$ws.Range('B741').Value = 'Nochistlán De Mejía'

# Remove trailing metadata rows (755-759); row 753 remains the last data row
$ws.Range("A755:D759").EntireRow.Delete()
